# Add season-record columns (Wins/Losses/Ties) to the stats table.
#
# The previous version of this workbook only pulled down team statistics;
# this adds the team's season record (Wins, Losses, Ties) as three new
# trailing columns, matching it onto every player row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, appended right after the existing "Unnamed: 28" column (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same look (bold, centered, bordered) as the rest
# of the header row by copying the formatting from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record for this team: 68 wins, 94 losses, 0 ties.
# Stamp it onto every player row (2 through 60).
$wins = 68
$losses = 94
$ties = 0

$lastRow = 60
for ($r = 2; $r -le $lastRow; $r++) {
  $ws.Cells.Item($r, 30).Value = $wins
  $ws.Cells.Item($r, 31).Value = $losses
  $ws.Cells.Item($r, 32).Value = $ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
